# Daily attendance processing - 2026-01-15 12:56:55
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session's attendance, as a comma separated string. The automated
# attendance processor normalizes the ordering of that list (moving the
# "System" account to the position the log replay produced) without
# touching any other column.
#
# This re-applies that normalization: every "Recorded By" cell whose text
# is exactly one of the two known stale orderings is rewritten to the
# corrected ordering. Cells that already have a different combination of
# authors (e.g. "admin@admin.com, System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$map = @{
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
